$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update raw measurement (B) and meter (C) values for trials 1-5.
# Column D holds formulas (=C/100) and will recalculate automatically.
$ws.Range("B2").Value = 33.01
$ws.Range("C2").Value = 3300

$ws.Range("B3").Value = 33
$ws.Range("C3").Value = 3303

$ws.Range("B4").Value = 33.03
$ws.Range("C4").Value = 3300

$ws.Range("B5").Value = 32.98
$ws.Range("C5").Value = 3299

$ws.Range("B6").Value = 33
$ws.Range("C6").Value = 3293

# Update the selected cell shown in the saved sheet view.
$ws.Range("E3").Select()
